# Show Advancing Fire Markers in Status Bar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Update the descriptive text in B34 (event e033 - Placing Advancing Fire Markers)
$newText = "<Bold>e033 Placing Advancing Fire Markers</Bold> `n<InlineUIContainer><Button Content='r4.61' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nPlace Advancing Fire Markers available to you per `n<InlineUIContainer><Button Content='r22.12' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. `nClick one of highlighted regions to place. `n<LineBreak/><LineBreak/>`nYou place up to six minus one marker for every three friendly tank losses (rounded up) . You may place more than one in a zone. `n<LineBreak/><LineBreak/>`nThe status bar on the bottom shows how many are remaining to place."

$ws.Range("B34").Value = $newText

# Increase row 34 height to fit the expanded text
$ws.Rows.Item(34).RowHeight = 150

# Update the active selection shown when the workbook is opened
$ws.Range("B35").Select()
